$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph entirely (it sat right after the
#    title heading).
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph with the page title text right before the
#    trailing "Feature image prompt for DALLE" paragraph.
$n = $d.Paragraphs.Count
$dalle = $d.Paragraphs.Item($n)
$insertionPoint = $d.Range($dalle.Range.Start, $dalle.Range.Start)
$xmlFrag = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play 'All Lucky Clover' Free - A Simple Yet Rewarding Slot Game</w:t></w:r></w:p><w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>"
$insertionPoint.InsertXML($xmlFrag)

# InsertXML splits the paragraph mark into two paragraphs; remove the stray
# empty paragraph it leaves behind between the new text and the DALLE prompt.
$strayIndex = $d.Paragraphs.Count - 1
$stray = $d.Paragraphs.Item($strayIndex)
$stray.Range.Delete()

# 3. Replace the DALLE image-prompt text with the meta-description copy,
#    keeping the run's existing italic formatting intact.
$dalle = $d.Paragraphs.Item($d.Paragraphs.Count)
$dalleText = $d.Range($dalle.Range.Start, $dalle.Range.End - 1)
$dalleText.Text = "Discover the simplistic yet rewarding gameplay of 'All Lucky Clover', featuring selectable paylines, Wild symbols, Scatters and a Gamble function. Play now for free."

Write-Output "done"
